$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" '29.159.42'
$ws.Range("E2").Value = '  -0.17%  '
Set-TextCell "D3" '1.838.12'
$ws.Range("E3").Value = '  -0.53%  '
Set-TextCell "D4" '0.9990'
$ws.Range("E4").Value = '  -0.02%  '
Set-TextCell "D5" '241.15'
$ws.Range("E5").Value = '  -2.12%  '
Set-TextCell "D6" '0.6866'
$ws.Range("E6").Value = '  -1.78%  '
Set-TextCell "D7" '0.9996'
$ws.Range("E7").Value = '  -0.02%  '
Set-TextCell "D8" '0.3016'
$ws.Range("E8").Value = '  -1.43%  '
Set-TextCell "D9" '0.07482'
$ws.Range("E9").Value = '  -3.15%  '
Set-TextCell "D10" '23.17'
$ws.Range("E10").Value = '  -1.62%  '
Set-TextCell "D11" '0.07666'
$ws.Range("E11").Value = '  -2.06%  '
Set-TextCell "D12" '1.838.53'
$ws.Range("E12").Value = '  -0.51%  '
Set-TextCell "D13" '5.062'
$ws.Range("E13").Value = '  -1.31%  '
Set-TextCell "D14" '0.6828'
$ws.Range("E14").Value = '  -0.57%  '
Set-TextCell "D15" '87.72'
$ws.Range("E15").Value = '  -5.99%  '
Set-TextCell "D16" '6.161'
$ws.Range("E16").Value = '  -7.15%  '
Set-TextCell "D17" '29.146.29'
$ws.Range("E17").Value = '  -0.21%  '
Set-TextCell "D18" '0.000008179'
$ws.Range("E18").Value = '  -1.72%  '
Set-TextCell "D19" '2.081.09'
$ws.Range("E19").Value = '  -0.45%  '
Set-TextCell "D20" '228.15'
$ws.Range("E20").Value = '  -5.47%  '
Set-TextCell "D21" '12.54'
$ws.Range("E21").Value = '  -1.79%  '
Set-TextCell "D22" '1.0000'
$ws.Range("E22").Value = '  +0.04%  '
Set-TextCell "D23" '7.405'
$ws.Range("E23").Value = '  -1.53%  '
Set-TextCell "D24" '0.9991'
$ws.Range("E24").Value = '  -0.07%  '
Set-TextCell "D25" '0.1456'
$ws.Range("E25").Value = '  -3.96%  '
Set-TextCell "D26" '159.87'
$ws.Range("E26").Value = '  +0.50%  '
Set-TextCell "D27" '8.771'
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("E28").Value = '  -1.08%  '
Set-TextCell "D29" '1.513'
$ws.Range("E29").Value = '  -2.04%  '
Set-TextCell "D30" '4.278'
$ws.Range("E30").Value = '  +1.00%  '
Set-TextCell "D31" '4.154'
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("E32").Value = '  +0.49%  '
Set-TextCell "D33" '0.05187'
$ws.Range("E33").Value = '  +1.21%  '
Set-TextCell "D34" '0.7656'
Set-TextCell "D35" '1.847'
$ws.Range("E35").Value = '  -1.35%  '
Set-TextCell "D36" '1.135'
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("E37").Value = '  -0.69%  '
Set-TextCell "D38" '1.317.93'
$ws.Range("E38").Value = '  +0.12%  '
Set-TextCell "D39" '0.01836'
$ws.Range("E39").Value = '  -1.98%  '
Set-TextCell "D40" '2.723'
$ws.Range("E40").Value = '  +0.42%  '
Set-TextCell "D41" '0.9331'
$ws.Range("E41").Value = '  -1.62%  '
Set-TextCell "D42" '5.798'
$ws.Range("E42").Value = '  -4.50%  '
Set-TextCell "D43" '104.68'
$ws.Range("E43").Value = '  -2.82%  '
Set-TextCell "D44" '0.9998'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell "D46" '65.19'
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell "D47" '1.982.99'
$ws.Range("E47").Value = '  -0.15%  '
Set-TextCell "D48" '0.5192'
$ws.Range("E48").Value = '  +0.31%  '
Set-TextCell "D49" '9.557'
$ws.Range("E49").Value = '  -1.78%  '
Set-TextCell "D50" '1.772'
$ws.Range("E50").Value = '  +0.31%  '
Set-TextCell "D51" '0.05938'
$ws.Range("E51").Value = '  +0.97%  '
